# Periodic fever syndromes.xlsx — add a "metadata" tab with the PanelApp
# query provenance, and refresh the "time_taken" (col F) timestamps on the
# "data" tab to the re-query time recorded in that new tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the F-column ("time_taken") timestamps on the "data" sheet.
# ---------------------------------------------------------------------
$timeTaken = @{
    2 = "2021-10-05 14:22:11.392598"
    3 = "2021-10-05 14:22:11.392606"
    4 = "2021-10-05 14:22:11.392609"
    5 = "2021-10-05 14:22:11.392612"
    6 = "2021-10-05 14:22:11.392615"
    7 = "2021-10-05 14:22:11.392617"
    8 = "2021-10-05 14:22:11.392620"
    9 = "2021-10-05 14:22:11.392622"
    10 = "2021-10-05 14:22:11.392625"
    11 = "2021-10-05 14:22:11.392628"
    12 = "2021-10-05 14:22:11.392630"
    13 = "2021-10-05 14:22:11.392633"
    14 = "2021-10-05 14:22:11.392635"
    15 = "2021-10-05 14:22:11.392638"
    16 = "2021-10-05 14:22:11.392640"
    17 = "2021-10-05 14:22:11.392643"
    18 = "2021-10-05 14:22:11.392646"
    19 = "2021-10-05 14:22:11.392648"
    20 = "2021-10-05 14:22:11.392651"
    21 = "2021-10-05 14:22:11.392653"
    22 = "2021-10-05 14:22:11.392655"
    23 = "2021-10-05 14:22:11.392658"
    24 = "2021-10-05 14:22:11.392661"
    25 = "2021-10-05 14:22:11.392663"
    26 = "2021-10-05 14:22:11.392666"
    27 = "2021-10-05 14:22:11.392669"
    28 = "2021-10-05 14:22:11.392671"
    29 = "2021-10-05 14:22:11.392673"
    30 = "2021-10-05 14:22:11.392676"
    31 = "2021-10-05 14:22:11.392678"
}

foreach ($row in $timeTaken.Keys) {
    $ws1.Cells.Item($row, 6).Value = $timeTaken[$row]
}

# ---------------------------------------------------------------------
# 2. Add the new "metadata" worksheet, after the existing "data" sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "metadata"

# Header row (B1:G1)
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Index column value (A2), same convention as the "data" sheet's column A
$ws2.Cells.Item(2, 1).Value = 0

# Match the header/index cell styling (bold, centered, top-aligned, boxed)
# used on the "data" sheet by copying its format onto the new cells.
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (row 2)
$ws2.Cells.Item(2, 2).Value = "Periodic fever syndromes"
$ws2.Cells.Item(2, 3).Value = 60

# data_version is stored as text ("1.15"), not a number, on the source sheet.
$dataVersion = $ws2.Cells.Item(2, 4)
$dataVersion.NumberFormat = "@"
$dataVersion.Value = "1.15"
$dataVersion.Style = "Normal"

$ws2.Cells.Item(2, 5).Value = "2021-08-03T17:12:09.143605Z"
$ws2.Cells.Item(2, 6).Value = "2021-10-05 14:22:11.388952"
$ws2.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/60/?format=json"

# Leave the "data" tab selected/active, matching the original workbook.
$ws1.Select()
